$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = 458
$ws.Range("D13").Value = 987
$ws.Range("D26").Value = 1547
$ws.Range("D38").Value = 1797
$ws.Range("D63").Value = 1989

$ws.Range("P21").Select()
